$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns are treated as text so numeric-looking strings
# (e.g. "528.61", "1.00") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '58.035.08'
$ws.Range("E2").Value = '  +2.34%  '
# Row 3
$ws.Range("D3").Value = '3.065.17'
$ws.Range("E3").Value = '  +2.64%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
$ws.Range("D5").Value = '528.61'
$ws.Range("E5").Value = '  +6.36%  '
# Row 6
$ws.Range("D6").Value = '143.78'
$ws.Range("E6").Value = '  +6.51%  '
# Row 7
$ws.Range("E7").Value = '  +0.01%  '
# Row 8
$ws.Range("E8").Value = '  +5.35%  '
# Row 9
$ws.Range("D9").Value = '7.65'
$ws.Range("E9").Value = '  +5.83%  '
# Row 10
$ws.Range("E10").Value = '  +7.51%  '
# Row 11
$ws.Range("D11").Value = '0.372'
$ws.Range("E11").Value = '  +6.08%  '
# Row 12
$ws.Range("E12").Value = '  +2.05%  '
# Row 13
$ws.Range("D13").Value = '3.588.63'
$ws.Range("E13").Value = '  +2.72%  '
# Row 14
$ws.Range("D14").Value = '27.47'
$ws.Range("E14").Value = '  +8.31%  '
# Row 15
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  +16.24%  '
# Row 16
$ws.Range("D16").Value = '58.018.23'
$ws.Range("E16").Value = '  +2.46%  '
# Row 17
$ws.Range("D17").Value = '6.23'
$ws.Range("E17").Value = '  +8.61%  '
# Row 18
$ws.Range("D18").Value = '3.082.21'
$ws.Range("E18").Value = '  +3.41%  '
# Row 19
$ws.Range("D19").Value = '13.24'
$ws.Range("E19").Value = '  +7.23%  '
# Row 20
$ws.Range("D20").Value = '8.25'
$ws.Range("E20").Value = '  +5.75%  '
# Row 21
$ws.Range("D21").Value = '342.65'
$ws.Range("E21").Value = '  +5.04%  '
# Row 22
$ws.Range("E22").Value = '  +0.16%  '
# Row 23
$ws.Range("E23").Value = '  -1.18%  '
# Row 24
$ws.Range("D24").Value = '0.506'
$ws.Range("E24").Value = '  +7.62%  '
# Row 25
$ws.Range("D25").Value = '65.49'
$ws.Range("E25").Value = '  +6.04%  '
# Row 26
$ws.Range("D26").Value = '0.0₃0985'
$ws.Range("E26").Value = '  +9.61%  '
# Row 27
$ws.Range("E27").Value = '  +4.98%  '
# Row 28
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.16%  '
# Row 29
$ws.Range("D29").Value = '7.04'
$ws.Range("E29").Value = '  +8.91%  '
# Row 30
$ws.Range("D30").Value = '7.45'
$ws.Range("E30").Value = '  +10.14%  '
# Row 31
$ws.Range("E31").Value = '  +7.44%  '
# Row 32
$ws.Range("E32").Value = '  +6.35%  '
# Row 33
$ws.Range("D33").Value = '21.23'
$ws.Range("E33").Value = '  +3.40%  '
# Row 34
$ws.Range("D34").Value = '4.82'
$ws.Range("E34").Value = '  +8.26%  '
# Row 35
$ws.Range("D35").Value = '157.34'
$ws.Range("E35").Value = '  +3.52%  '
# Row 36
$ws.Range("D36").Value = '6.00'
$ws.Range("E36").Value = '  +7.08%  '
# Row 37
$ws.Range("D37").Value = '1.34'
$ws.Range("E37").Value = '  +4.35%  '
# Row 38
$ws.Range("D38").Value = '26.28'
$ws.Range("E38").Value = '  +12.44%  '
# Row 39
$ws.Range("E39").Value = '  +4.81%  '
# Row 40
$ws.Range("D40").Value = '3.100.03'
$ws.Range("E40").Value = '  +2.77%  '
# Row 41
$ws.Range("D41").Value = '37.78'
$ws.Range("E41").Value = '  +3.56%  '
# Row 42
$ws.Range("D42").Value = '3.97'
$ws.Range("E42").Value = '  +11.43%  '
# Row 43
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.48'
$ws.Range("E43").Value = '  +5.84%  '
# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.348.69'
$ws.Range("E44").Value = '  +5.64%  '
# Row 45
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '1.05'
$ws.Range("E45").Value = '  +4.47%  '
# Row 46
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.09%  '
# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.667'
$ws.Range("E47").Value = '  +4.05%  '
# Row 48
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '6.15'
$ws.Range("E48").Value = '  +7.41%  '
# Row 49
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '2.02'
$ws.Range("E49").Value = '  +3.49%  '
# Row 50
$ws.Range("E50").Value = '  +4.44%  '
# Row 51
$ws.Range("D51").Value = '20.41'
$ws.Range("E51").Value = '  +7.29%  '

# Restore default (Normal) style so no stray number-format style
# is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
